$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# xlCenter = -4108 (horizontal center alignment), matching the existing
# header/body rows on the sheet (style index 1 in styles.xml).
$xlCenter = -4108

# --- New "Coeff roue vitesse" block (rows 10-12) -------------------------
$ws.Range("A10").Value = "Coeff roue vitesse"

$ws.Range("A11").Value = "P"
$ws.Range("B11").Value = "D"
$ws.Range("C11").Value = "I"
$ws.Range("A11:C11").HorizontalAlignment = $xlCenter

$ws.Range("A12").Value = "'2.5 "
$ws.Range("B12").Value = "'0.01"
$ws.Range("C12").Value = "'0.05"
$ws.Range("D12").Value = ""
$ws.Range("E12").Value = "Nerveux"
$ws.Range("A12:E12").Style = "Normal"
$ws.Range("A12:E12").HorizontalAlignment = $xlCenter

# Empty, centered placeholder cells left below the new block.
$ws.Range("B13").Value = ""
$ws.Range("B13").HorizontalAlignment = $xlCenter

$ws.Range("B15").Value = ""
$ws.Range("B15").HorizontalAlignment = $xlCenter

$ws.Range("B16").Value = ""
$ws.Range("B16").HorizontalAlignment = $xlCenter

$ws.Range("B17").Value = ""
$ws.Range("B17").HorizontalAlignment = $xlCenter

# --- Column A width (author widened/auto-fit it for the new text) -------
$ws.Columns.Item(1).ColumnWidth = 12.8

# --- Selection matches the author's saved cursor position ---------------
$ws.Range("E11").Select()
